$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.023065755192931
$ws.Range("C2").Value = 0.1927725476463422
$ws.Range("D2").Value = 0.207869663336183
$ws.Range("E2").Value = 0.1706476245223243
$ws.Range("F2").Value = 1.255145777563058
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("J2").Value = 0.1839599612594327
$ws.Range("N2").Value = 1.056347490988195
$ws.Range("O2").Value = 2.887423348435988

$ws.Range("B3").Value = 0.9226797810767948
$ws.Range("C3").Value = 0.168586055212387
$ws.Range("D3").Value = 0.2033822855059952
$ws.Range("E3").Value = 0.1668010292406663
$ws.Range("F3").Value = 1.248864296630515
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("J3").Value = 0.1796355010563033
$ws.Range("N3").Value = 1.06427707386684
$ws.Range("O3").Value = 2.881533560150899

$ws.Range("B4").Value = 0.8611533150814239
$ws.Range("C4").Value = 0.1536945318088669
$ws.Range("D4").Value = 0.2007029652277197
$ws.Range("E4").Value = 0.1645237336686982
$ws.Range("F4").Value = 1.245831171726408
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("J4").Value = 0.1770892791536411
$ws.Range("N4").Value = 1.069608717380113
$ws.Range("O4").Value = 2.879976776102922

$ws.Range("B5").Value = 0.8361099128395324
$ws.Range("C5").Value = 0.1476161486325509
$ws.Range("D5").Value = 0.1996303047372407
$ws.Range("E5").Value = 0.1636169988557654
$ws.Range("F5").Value = 1.244802045488129
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("J5").Value = 0.1760790637621525
$ws.Range("N5").Value = 1.071897886162361
$ws.Range("O5").Value = 2.879859350500169

$ws.Range("B6").Value = 0.83195327108524
$ws.Range("C6").Value = 0.14660624532263
$ws.Range("D6").Value = 0.1994533515173202
$ws.Range("E6").Value = 0.163467722363432
$ws.Range("F6").Value = 1.244643649030849
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("J6").Value = 0.1759129719723163
$ws.Range("N6").Value = 1.072285038918089
$ws.Range("O6").Value = 2.879871049352971

$ws.Range("B7").Value = 0.8608154514010948
$ws.Range("C7").Value = 0.1536125965312181
$ws.Range("D7").Value = 0.2006884211611464
$ws.Range("E7").Value = 0.1645114189223129
$ws.Range("F7").Value = 1.245816455169148
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("J7").Value = 0.1770755441597842
$ws.Range("N7").Value = 1.069639118158598
$ws.Range("O7").Value = 2.879973100468391

$ws.Range("B8").Value = 0.9884303583446581
$ws.Range("C8").Value = 0.1844417240303926
$ws.Range("D8").Value = 0.2063067025761427
$ws.Range("E8").Value = 0.1693037883256707
$ws.Range("F8").Value = 1.252808809370478
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("J8").Value = 0.1824462465183885
$ws.Range("N8").Value = 1.058985601662734
$ws.Range("O8").Value = 2.884964521375338

$ws.Range("B9").Value = 1.239521282611065
$ws.Range("C9").Value = 0.2445615781981303
$ws.Range("D9").Value = 0.217923678322606
$ws.Range("E9").Value = 0.1793718603239682
$ws.Range("F9").Value = 1.273070072568629
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("J9").Value = 0.19384493790146
$ws.Range("N9").Value = 1.041762886983058
$ws.Range("O9").Value = 2.911141077063661

$ws.Range("B10").Value = 1.424471702276662
$ws.Range("C10").Value = 0.2885159083516555
$ws.Range("D10").Value = 0.2268211270900196
$ws.Range("E10").Value = 0.1871778862982154
$ws.Range("F10").Value = 1.291970975281259
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("J10").Value = 0.2027515769076587
$ws.Range("N10").Value = 1.031341643504327
$ws.Range("O10").Value = 2.940435813039301

$ws.Range("B11").Value = 1.508706660821929
$ws.Range("C11").Value = 0.3084629902243137
$ws.Range("D11").Value = 0.2309469371568298
$ws.Range("E11").Value = 0.1908180179119725
$ws.Range("F11").Value = 1.301446358644839
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("J11").Value = 0.206919811752897
$ws.Range("N11").Value = 1.027084709268799
$ws.Range("O11").Value = 2.955963794899986

$ws.Range("B12").Value = 1.540617658664985
$ws.Range("C12").Value = 0.3160092671066081
$ws.Range("D12").Value = 0.2325204645420342
$ws.Range("E12").Value = 0.1922092487579263
$ws.Range("F12").Value = 1.305160925925207
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("J12").Value = 0.2085150194369589
$ws.Range("N12").Value = 1.025542229027337
$ws.Range("O12").Value = 2.962161591097669

$ws.Range("B13").Value = 1.533744495399276
$ws.Range("C13").Value = 0.3143843694171551
$ws.Range("D13").Value = 0.2321810817753232
$ws.Range("E13").Value = 0.1919090537716883
$ws.Range("F13").Value = 1.304355299006914
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("J13").Value = 0.208170715865819
$ws.Range("N13").Value = 1.025871337917167
$ws.Range("O13").Value = 2.960812638224326

$ws.Range("B14").Value = 1.511331744020652
$ws.Range("C14").Value = 0.309083974285187
$ws.Range("D14").Value = 0.2310761687946155
$ws.Range("E14").Value = 0.1909322190827609
$ws.Range("F14").Value = 1.301749422487248
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("J14").Value = 0.2070507138492985
$ws.Range("N14").Value = 1.026956415207792
$ws.Range("O14").Value = 2.956467318163533

$ws.Range("B15").Value = 1.497604945947444
$ws.Range("C15").Value = 0.3058363745898305
$ws.Range("D15").Value = 0.2304008300599634
$ws.Range("E15").Value = 0.1903355447879846
$ws.Range("F15").Value = 1.300169723547185
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("J15").Value = 0.2063668673020942
$ws.Range("N15").Value = 1.027630110142908
$ws.Range("O15").Value = 2.953847088695028

$ws.Range("B16").Value = 1.418968648510713
$ws.Range("C16").Value = 0.2872113242821968
$ws.Range("D16").Value = 0.2265530639959934
$ws.Range("E16").Value = 0.186941786549113
$ws.Range("F16").Value = 1.291369414790452
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("J16").Value = 0.2024815207011841
$ws.Range("N16").Value = 1.031629575323279
$ws.Range("O16").Value = 2.939465419121149

$ws.Range("B17").Value = 1.370752502622508
$ws.Range("C17").Value = 0.2757729229547294
$ws.Range("D17").Value = 0.2242125792819252
$ws.Range("E17").Value = 0.1848826366529153
$ws.Range("F17").Value = 1.286195609664333
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("J17").Value = 0.2001278590133779
$ws.Range("N17").Value = 1.034206990098951
$ws.Range("O17").Value = 2.931207404103418

$ws.Range("B18").Value = 1.34302935130313
$ws.Range("C18").Value = 0.2691893635192173
$ws.Range("D18").Value = 0.2228737711623978
$ws.Range("E18").Value = 0.1837066590381795
$ws.Range("F18").Value = 1.283302328722229
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("J18").Value = 0.1987850646699343
$ws.Range("N18").Value = 1.035734984557457
$ws.Range("O18").Value = 2.926664777611506

$ws.Range("B19").Value = 1.333644441410115
$ws.Range("C19").Value = 0.2669595224290333
$ws.Range("D19").Value = 0.2224217434241069
$ws.Range("E19").Value = 0.1833099351380199
$ws.Range("F19").Value = 1.282336883414388
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("J19").Value = 0.1983323008137035
$ws.Range("N19").Value = 1.03626015878384
$ws.Range("O19").Value = 2.925162269982366

$ws.Range("B20").Value = 1.375884220152727
$ws.Range("C20").Value = 0.2769910279275223
$ws.Range("D20").Value = 0.2244609650526996
$ws.Range("E20").Value = 0.1851009684252247
$ws.Range("F20").Value = 1.286737823787703
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("J20").Value = 0.2003772748145281
$ws.Range("N20").Value = 1.033927907589678
$ws.Range("O20").Value = 2.932065035159724

$ws.Range("B21").Value = 1.517914568988601
$ws.Range("C21").Value = 0.3106410283907053
$ws.Range("D21").Value = 0.2314004060059176
$ws.Range("E21").Value = 0.1912187921367376
$ws.Range("F21").Value = 1.302511397546013
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("J21").Value = 0.2073792297691739
$ws.Range("N21").Value = 1.026635814959583
$ws.Range("O21").Value = 2.957735013381637

$ws.Range("B22").Value = 1.610815040999739
$ws.Range("C22").Value = 0.3325908269223987
$ws.Range("D22").Value = 0.2360008134489817
$ws.Range("E22").Value = 0.1952916936768858
$ws.Range("F22").Value = 1.313557496809636
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("J22").Value = 0.2120532754533002
$ws.Range("N22").Value = 1.022275264201497
$ws.Range("O22").Value = 2.97636405478417

$ws.Range("B23").Value = 1.561225851994607
$ws.Range("C23").Value = 0.3208798042342664
$ws.Range("D23").Value = 0.2335395638139346
$ws.Range("E23").Value = 0.1931110963544427
$ws.Range("F23").Value = 1.307594438182406
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("J23").Value = 0.209549686055496
$ws.Range("N23").Value = 1.024565500921149
$ws.Range("O23").Value = 2.96625154534496

$ws.Range("B24").Value = 1.373564180044809
$ws.Range("C24").Value = 0.2764403459136702
$ws.Range("D24").Value = 0.2243486487598148
$ws.Range("E24").Value = 0.1850022361624823
$ws.Range("F24").Value = 1.286492435841666
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("J24").Value = 0.2002644816727752
$ws.Range("N24").Value = 1.034053936861433
$ws.Range("O24").Value = 2.931676661615342

$ws.Range("B25").Value = 1.171508640567311
$ws.Range("C25").Value = 0.2283346660701966
$ws.Range("D25").Value = 0.2147171087037236
$ws.Range("E25").Value = 0.17657638314369
$ws.Range("F25").Value = 1.2668854124083
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("J25").Value = 0.1906681233375451
$ws.Range("N25").Value = 1.046029786764343
$ws.Range("O25").Value = 2.902297696863002
